$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 453.2580555555555
$ws.Cells.Item(2, 4).Value = 33.08451422318829
$ws.Cells.Item(2, 5).Value = 390
$ws.Cells.Item(2, 6).Value = 430
$ws.Cells.Item(2, 7).Value = 448
$ws.Cells.Item(2, 8).Value = 470
$ws.Cells.Item(2, 9).Value = 733

$ws.Cells.Item(3, 3).Value = 35.78553560606061
$ws.Cells.Item(3, 4).Value = 4.437912724794431
$ws.Cells.Item(3, 5).Value = 22.73
$ws.Cells.Item(3, 6).Value = 32.36
$ws.Cells.Item(3, 7).Value = 35.78
$ws.Cells.Item(3, 8).Value = 39.07
$ws.Cells.Item(3, 9).Value = 50.85

$ws.Cells.Item(4, 3).Value = 1.750667171717172
$ws.Cells.Item(4, 4).Value = 1.101502338716308
$ws.Cells.Item(4, 5).Value = 0.09
$ws.Cells.Item(4, 6).Value = 0.97
$ws.Cells.Item(4, 7).Value = 1.45
$ws.Cells.Item(4, 8).Value = 2.31
$ws.Cells.Item(4, 9).Value = 10.39

$ws.Cells.Item(5, 3).Value = 320.2051017676768
$ws.Cells.Item(5, 4).Value = 5.603780927290645
$ws.Cells.Item(5, 5).Value = 308.26
$ws.Cells.Item(5, 6).Value = 316.52
$ws.Cells.Item(5, 7).Value = 319.9
$ws.Cells.Item(5, 8).Value = 325.23
$ws.Cells.Item(5, 9).Value = 332.88

$ws.Cells.Item(6, 3).Value = 28.71858712121212
$ws.Cells.Item(6, 4).Value = 1.995589405843883
$ws.Cells.Item(6, 5).Value = 22.33
$ws.Cells.Item(6, 6).Value = 27.09
$ws.Cells.Item(6, 7).Value = 28.82
$ws.Cells.Item(6, 8).Value = 29.82
$ws.Cells.Item(6, 9).Value = 42.06

$ws.Cells.Item(7, 3).Value = -70.23260101010101
$ws.Cells.Item(7, 4).Value = 19.48500092970017
$ws.Cells.Item(7, 5).Value = -121
$ws.Cells.Item(7, 6).Value = -86
$ws.Cells.Item(7, 7).Value = -69
$ws.Cells.Item(7, 8).Value = -50
$ws.Cells.Item(7, 9).Value = -35

$ws.Cells.Item(8, 3).Value = 9.808034157802984
$ws.Cells.Item(8, 4).Value = 3.354816805637959
$ws.Cells.Item(8, 5).Value = -19
$ws.Cells.Item(8, 6).Value = 8.800000000000001
$ws.Cells.Item(8, 7).Value = 10.2
$ws.Cells.Item(8, 8).Value = 11.8
$ws.Cells.Item(8, 9).Value = 17

$ws.Cells.Item(9, 3).Value = 8.50320707070707
$ws.Cells.Item(9, 4).Value = 1.117449036735686
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 8
$ws.Cells.Item(9, 7).Value = 9
$ws.Cells.Item(9, 8).Value = 10
$ws.Cells.Item(9, 9).Value = 11

$ws.Cells.Item(10, 3).Value = 867.8014974747476
$ws.Cells.Item(10, 4).Value = 0.4588582799635839
$ws.Cells.Item(10, 5).Value = 867.1
$ws.Cells.Item(10, 6).Value = 867.5
$ws.Cells.Item(10, 7).Value = 867.9
$ws.Cells.Item(10, 8).Value = 868.3
$ws.Cells.Item(10, 9).Value = 868.5

$ws.Cells.Item(11, 3).Value = 0.2266959159595959
$ws.Cells.Item(11, 4).Value = 0.145131500830238
$ws.Cells.Item(11, 5).Value = 0.071936
$ws.Cells.Item(11, 6).Value = 0.133632
$ws.Cells.Item(11, 7).Value = 0.246784
$ws.Cells.Item(11, 8).Value = 0.452608
$ws.Cells.Item(11, 9).Value = 0.987136

$ws.Cells.Item(12, 3).Value = 22.74816919191919
$ws.Cells.Item(12, 4).Value = 12.36859170704831
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = 10
$ws.Cells.Item(12, 7).Value = 23
$ws.Cells.Item(12, 8).Value = 37
$ws.Cells.Item(12, 9).Value = 40

$ws.Cells.Item(13, 3).Value = 0.661060606060606
$ws.Cells.Item(13, 4).Value = 0.7506905444337185
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 1
$ws.Cells.Item(13, 9).Value = 2

$ws.Cells.Item(14, 3).Value = 1.830025252525252
$ws.Cells.Item(14, 4).Value = 1.690580924950423
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 8).Value = 2
$ws.Cells.Item(14, 9).Value = 5

$ws.Cells.Item(15, 3).Value = 87.492601010101
$ws.Cells.Item(15, 4).Value = 19.48500092970293
$ws.Cells.Item(15, 5).Value = 52.26
$ws.Cells.Item(15, 6).Value = 67.25999999999999
$ws.Cells.Item(15, 7).Value = 86.25999999999999
$ws.Cells.Item(15, 8).Value = 103.26
$ws.Cells.Item(15, 9).Value = 138.26

$ws.Cells.Item(16, 3).Value = -80.6473610253542
$ws.Cells.Item(16, 4).Value = 18.81304349493355
$ws.Cells.Item(16, 5).Value = -121.8240101588271
$ws.Cells.Item(16, 6).Value = -96.96183611348224
$ws.Cells.Item(16, 7).Value = -80.45410721860875
$ws.Cells.Item(16, 8).Value = -61.33779541063677
$ws.Cells.Item(16, 9).Value = -42.95746206410165

$ws.Cells.Item(17, 3).Value = -70.83932686755121
$ws.Cells.Item(17, 4).Value = 20.02536296935124
$ws.Cells.Item(17, 5).Value = -138.8714834061964
$ws.Cells.Item(17, 6).Value = -86.5149694202523
$ws.Cells.Item(17, 7).Value = -69.49305820175223
$ws.Cells.Item(17, 8).Value = -50.33195619988427
$ws.Cells.Item(17, 9).Value = -35.57382219273629

